$wb = $excel.ActiveWorkbook

# This script applies updated commodity/leve-profit price data to several
# sheets, as produced by the scheduled pricing-data runner.

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 5904.1113
$ws.Range("I61").Value = 784.25
$ws.Range("K61").Value = 2352.75
$ws.Range("M61").Value = -2180.75
$ws.Range("H107").Value = 1539.25
$ws.Range("I107").Value = 1719.4445
$ws.Range("K107").Value = 1719.4445
$ws.Range("M107").Value = 200.5554999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3268.8333
$ws.Range("J2").Value = 4706.5
$ws.Range("L2").Value = 4706.5
$ws.Range("N2").Value = -4932.5
$ws.Range("H32").Value = 4283.909
$ws.Range("I32").Value = 4283.909
$ws.Range("K32").Value = 4283.909
$ws.Range("M32").Value = -3996.909
$ws.Range("H62").Value = 90624.5
$ws.Range("J62").Value = 90624.5
$ws.Range("L62").Value = 90624.5
$ws.Range("N62").Value = -91872.5
$ws.Range("H63").Value = 12116.458
$ws.Range("I63").Value = 31427.285
$ws.Range("J63").Value = 4164.9414
$ws.Range("K63").Value = 31427.285
$ws.Range("L63").Value = 4164.9414
$ws.Range("M63").Value = -30741.285
$ws.Range("N63").Value = -5536.9414
$ws.Range("H65").Value = 90624.5
$ws.Range("J65").Value = 90624.5
$ws.Range("L65").Value = 271873.5
$ws.Range("N65").Value = -278113.5
$ws.Range("H66").Value = 12116.458
$ws.Range("I66").Value = 31427.285
$ws.Range("J66").Value = 4164.9414
$ws.Range("K66").Value = 157136.425
$ws.Range("L66").Value = 20824.707
$ws.Range("M66").Value = -153704.425
$ws.Range("N66").Value = -27688.707
$ws.Range("H74").Value = 1970.0312
$ws.Range("I74").Value = 2253.65
$ws.Range("K74").Value = 2253.65
$ws.Range("M74").Value = -1379.65
$ws.Range("H77").Value = 1970.0312
$ws.Range("I77").Value = 2253.65
$ws.Range("K77").Value = 11268.25
$ws.Range("M77").Value = -6900.25
$ws.Range("H116").Value = 3268.8333
$ws.Range("J116").Value = 4706.5
$ws.Range("L116").Value = 4706.5
$ws.Range("N116").Value = -9294.5
$ws.Range("H132").Value = 8999.666999999999
$ws.Range("I132").Value = 8999.666999999999
$ws.Range("K132").Value = 26999.001
$ws.Range("M132").Value = -24469.001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3268.8333
$ws.Range("J3").Value = 4706.5
$ws.Range("L3").Value = 4706.5
$ws.Range("N3").Value = -4934.5
$ws.Range("H82").Value = 20466.166
$ws.Range("J82").Value = 29949.666
$ws.Range("L82").Value = 29949.666
$ws.Range("N82").Value = -30715.666
$ws.Range("H85").Value = 20466.166
$ws.Range("J85").Value = 29949.666
$ws.Range("L85").Value = 29949.666
$ws.Range("N85").Value = -32601.666
$ws.Range("H86").Value = 7824.1665
$ws.Range("I86").Value = 1951.1818
$ws.Range("J86").Value = 12793.615
$ws.Range("K86").Value = 1951.1818
$ws.Range("L86").Value = 12793.615
$ws.Range("M86").Value = -828.1818000000001
$ws.Range("N86").Value = -15039.615
$ws.Range("H89").Value = 7824.1665
$ws.Range("I89").Value = 1951.1818
$ws.Range("J89").Value = 12793.615
$ws.Range("K89").Value = 9755.909
$ws.Range("L89").Value = 63968.075
$ws.Range("M89").Value = -4139.909
$ws.Range("N89").Value = -75200.075
$ws.Range("H94").Value = 1461.6522
$ws.Range("I94").Value = 1092.9286
$ws.Range("K94").Value = 1092.9286
$ws.Range("M94").Value = -641.9286
$ws.Range("H134").Value = 5510.4346
$ws.Range("I134").Value = 5510.4346
$ws.Range("K134").Value = 16531.3038
$ws.Range("M134").Value = -13996.3038

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 5000166
$ws.Range("I7").Value = 6579118.5
$ws.Range("K7").Value = 6579118.5
$ws.Range("M7").Value = -6579005.5
$ws.Range("H22").Value = 810.6
$ws.Range("I22").Value = 817
$ws.Range("J22").Value = 801
$ws.Range("K22").Value = 817
$ws.Range("L22").Value = 801
$ws.Range("M22").Value = -467
$ws.Range("N22").Value = -1501
$ws.Range("H58").Value = 2260.6
$ws.Range("I58").Value = 2289.5557
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 2289.5557
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -2086.5557
$ws.Range("N58").Value = -2406
$ws.Range("H69").Value = 14741.25
$ws.Range("I69").Value = 14482.5
$ws.Range("K69").Value = 14482.5
$ws.Range("M69").Value = -13733.5
$ws.Range("H72").Value = 14741.25
$ws.Range("I72").Value = 14482.5
$ws.Range("K72").Value = 43447.5
$ws.Range("M72").Value = -39703.5
$ws.Range("H99").Value = 1724.625
$ws.Range("I99").Value = 1611.8572
$ws.Range("K99").Value = 1611.8572
$ws.Range("M99").Value = -113.8571999999999
$ws.Range("H126").Value = 1724.625
$ws.Range("I126").Value = 1611.8572
$ws.Range("K126").Value = 4835.571599999999
$ws.Range("M126").Value = -2365.571599999999
$ws.Range("H136").Value = 2260.6
$ws.Range("I136").Value = 2289.5557
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 6868.6671
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -4318.6671
$ws.Range("N136").Value = -11100

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5218.7
$ws.Range("J34").Value = 6699.4287
$ws.Range("L34").Value = 20098.2861
$ws.Range("N34").Value = -20266.2861
$ws.Range("H55").Value = 2615.818
$ws.Range("J55").Value = 4333
$ws.Range("L55").Value = 12999
$ws.Range("N55").Value = -13353
$ws.Range("H129").Value = 2819.75
$ws.Range("I129").Value = 648
$ws.Range("J129").Value = 4991.5
$ws.Range("K129").Value = 1944
$ws.Range("L129").Value = 14974.5
$ws.Range("M129").Value = 3056
$ws.Range("N129").Value = -24974.5
$ws.Range("H136").Value = 4588.3335
$ws.Range("I136").Value = 3395
$ws.Range("K136").Value = 10185
$ws.Range("M136").Value = -5085

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1674.25
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1674.25
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 1674.25
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -2666.25
$ws.Range("H122").Value = 6369.3335
$ws.Range("I122").Value = 13500
$ws.Range("K122").Value = 40500
$ws.Range("M122").Value = -38050
$ws.Range("H126").Value = 8680.416999999999
$ws.Range("I126").Value = 3221.6667
$ws.Range("K126").Value = 9665.000100000001
$ws.Range("M126").Value = -7195.000100000001
$ws.Range("H132").Value = 2275.3635
$ws.Range("I132").Value = 1953.625
$ws.Range("J132").Value = 3133.3333
$ws.Range("K132").Value = 5860.875
$ws.Range("L132").Value = 9399.999899999999
$ws.Range("M132").Value = -3330.875
$ws.Range("N132").Value = -14459.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1801
$ws.Range("I7").Value = 1801.25
$ws.Range("K7").Value = 1801.25
$ws.Range("M7").Value = -1689.25
$ws.Range("H126").Value = 1801
$ws.Range("I126").Value = 1801.25
$ws.Range("K126").Value = 5403.75
$ws.Range("M126").Value = -2933.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2835.1
$ws.Range("I126").Value = 1521.1
$ws.Range("J126").Value = 4149.1
$ws.Range("K126").Value = 4563.299999999999
$ws.Range("L126").Value = 12447.3
$ws.Range("M126").Value = -2093.299999999999
$ws.Range("N126").Value = -17387.3
